# Switch to the first sheet ("8CH5Y") and make it the active tab
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("8CH5Y")
$ws1.Activate()

# Row 2 data corrections
$ws1.Range("A2").Value = 21671
$ws1.Range("B2").Value = "HELOISA PIASSALI DE CASTRO"
$ws1.Range("C2").Value = 5534467037
$ws1.Range("D2").Value = "LISLAINE PIASSALI DE CASTRO"

# Highlight the corrected CPF cell in yellow
$ws1.Range("C2").Interior.Color = 62207

# Selection moves to C2
$ws1.Range("C2").Select()

# Widen the new 4th column (same width as column B)
$ws1.Columns.Item(4).ColumnWidth = 53.8
